$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A (bank facility reference) - shared strings 19, 20
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("A6").Value = "Scone GEF"

# Column B (UKEF facility ID)
$ws.Range("B5").Value = 20001371
$ws.Range("B6").Value = 20001371

# Column C (exporter) - shared strings 21, 22
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("C6").Value = "Scone exporter"

# Column D (base currency)
$ws.Range("D5").Value = "GBP"
$ws.Range("D6").Value = "GBP"

# Column E (facility limit)
$ws.Range("E5").Value = 7000000
$ws.Range("E6").Value = 770000

# Column F (facility utilisation)
$ws.Range("F5").Value = 3938753.8
$ws.Range("F6").Value = 761579.37

# Column G (total fees accrued for the period)
$ws.Range("G5").Value = 777
$ws.Range("G6").Value = 777

# Column H (fees paid to UKEF for the period)
$ws.Range("H5").Value = 456
$ws.Range("H6").Value = 456.77

# Column I (fees paid to UKEF currency)
$ws.Range("I5").Value = "GBP"
$ws.Range("I6").Value = "GBP"

# Column J (payment currency)
$ws.Range("J5").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Apply matching row styles (same cellXfs as the data rows above, so row 5/6
# pick up the workbook's existing text/number formatting instead of General)
$ws.Range("A5:A6").Style = $ws.Range("A4").Style
$ws.Range("B5:B6").Style = $ws.Range("B4").Style
$ws.Range("C5:C6").Style = $ws.Range("C4").Style
$ws.Range("D5:D6").Style = $ws.Range("D4").Style
$ws.Range("E5:E6").Style = $ws.Range("E4").Style
$ws.Range("F5:F6").Style = $ws.Range("F4").Style
$ws.Range("G5:G6").Style = $ws.Range("G4").Style
$ws.Range("H5:H6").Style = $ws.Range("H4").Style
$ws.Range("I5:I6").Style = $ws.Range("I4").Style
$ws.Range("J5:J6").Style = $ws.Range("J4").Style

# Update selection to match the target view state
$ws.Range("F8").Select()
